$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.107.64'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.470.50'
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.83'
$ws.Range("E5").Value = '  +0.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.85'
$ws.Range("E6").Value = '  +0.90%  '
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.542'
$ws.Range("E8").Value = '  +0.67%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.470.48'
$ws.Range("E9").Value = '  +1.46%  '
$ws.Range("E10").Value = '  +1.11%  '
$ws.Range("E11").Value = '  +1.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.28'
$ws.Range("E12").Value = '  +0.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.354'
$ws.Range("E13").Value = '  +1.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.06'
$ws.Range("E14").Value = '  +8.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000178'
$ws.Range("E15").Value = '  -0.02%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.918.98'
$ws.Range("E16").Value = '  +2.28%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.121.16'
$ws.Range("E17").Value = '  +1.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.470.86'
$ws.Range("E18").Value = '  -1.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.92'
$ws.Range("E19").Value = '  -0.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.08'
$ws.Range("E20").Value = '  +2.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '329.85'
$ws.Range("E21").Value = '  +0.98%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.20'
$ws.Range("E23").Value = '  +8.05%  '
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.39'
$ws.Range("E25").Value = '  +1.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '660.09'
$ws.Range("E26").Value = '  +6.96%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.15'
$ws.Range("E27").Value = '  +8.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0996'
$ws.Range("E28").Value = '  +0.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.603.28'
$ws.Range("E29").Value = '  +3.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  +0.03%  '
$ws.Range("E31").Value = '  +3.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.12'
$ws.Range("E32").Value = '  +0.28%  '
$ws.Range("E33").Value = '  +2.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.136'
$ws.Range("E34").Value = '  +0.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.53'
$ws.Range("E35").Value = '  +3.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").Value = '  +0.15%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.79'
$ws.Range("E37").Value = '  +0.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.50'
$ws.Range("E38").Value = '  +1.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.373'
$ws.Range("E39").Value = '  +0.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '152.97'
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.81'
$ws.Range("E41").Value = '  +1.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.71'
$ws.Range("E42").Value = '  +1.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.76'
$ws.Range("E43").Value = '  +0.93%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0₆0300'
$ws.Range("E45").Value = '  +6.24%  '
$ws.Range("E46").Value = '  +27.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '147.29'
$ws.Range("E47").Value = '  +2.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.62'
$ws.Range("E48").Value = '  +1.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.80'
$ws.Range("E49").Value = '  +2.30%  '
$ws.Range("E50").Value = '  +1.41%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0516'
$ws.Range("E51").Value = '  +0.41%  '
